$d = $word.ActiveDocument

# 1. Update activation date 2012 -> 2024
$null = $d.Content.Find.Execute("Ativação: 01/01/2012", $false, $false, $false, $false, $false, $true, 1, $false, 'Ativação: 01/01/2024', 2)

# 2. Insert English translation paragraph after "Objetivos" body text
$pObjetivos = $null
foreach ($p in $d.Paragraphs) {
  if ($p.Range.Text.StartsWith("Introduzir o aluno")) { $pObjetivos = $p }
}
$null = $pObjetivos.Range.InsertParagraphAfter()
$pObjetivosEn = $pObjetivos.Next()
$objetivosEnText = 'Introduction to Chemical Reaction Engineering through the fundamental concepts of chemical kinetics applied to ideal chemical reactors'
$pObjetivosEn.Range.InsertBefore($objetivosEnText)
$objEnStart = $pObjetivosEn.Range.Start
$objEnRange = $d.Range($objEnStart, $objEnStart + $objetivosEnText.Length)
$objEnRange.Font.Italic = $true

# 3. Add new "Docente(s) Responsável(eis)" bullet line before the existing one
$pDocente = $null
foreach ($p in $d.Paragraphs) {
  if ($p.Range.Text.StartsWith("6310316")) { $pDocente = $p }
}
$docenteNewText = '5963230 - Leandro Gonçalves de Aguiar' + [char]11
$docenteStart = $pDocente.Range.Start
$docenteInsertRange = $d.Range($docenteStart, $docenteStart)
$docenteInsertRange.InsertBefore($docenteNewText)

# 4. Insert English translation paragraph after "Programa resumido" body text
$pResumido = $null
foreach ($p in $d.Paragraphs) {
  if ($p.Range.Text.StartsWith("1. Introdução a cinética")) { $pResumido = $p }
}
$null = $pResumido.Range.InsertParagraphAfter()
$pResumidoEn = $pResumido.Next()
$resumidoEnText = '1. Introduction to Kinetics. 2. Reactions at Constant Volume. 3. Reactions at Variable Volume. 4. Ideal Models of Isothermal Chemical Reactors. 5. Analysis of Kinetic Data in Isothermal Chemical Reactors'
$pResumidoEn.Range.InsertBefore($resumidoEnText)
$resEnStart = $pResumidoEn.Range.Start
$resEnRange = $d.Range($resEnStart, $resEnStart + $resumidoEnText.Length)
$resEnRange.Font.Italic = $true

# 5. Rewrite "Programa" body text (merge into one run, drop hour annotations) and add EN translation
$pPrograma = $null
foreach ($p in $d.Paragraphs) {
  if ($p.Range.Text.StartsWith("1. INTRODUÇÃO A CINÉTICA")) { $pPrograma = $p }
}
$programaPtText = '1. INTRODUÇÃO A CINÉTICATipos de Reações Químicas. Lei de velocidade e seus principais parâmetros. Influência da temperatura sobre a taxa da reação. Ativação das reações químicas Equação de Arrhenius. Energia de ativação. Conversão. Concentração e sua variação numa transformação química. 2. REAÇÕES A VOLUME CONSTANTEReações irreversíveis de ordem um. Reações irreversíveis de ordem dois. Reações irreversíveis de ordem três. Reações irreversíveis de ordem qualquer. 3. REAÇÕES A VOLUME VARIÁVELConceitos. Fração de conversão volumétrica. Reações a volume variável de ordem um e dois. 4. MODELOS IDEAIS DE REATORES QUÍMICOS ISOTÉRMICOS: Equações fundamentais de projeto de reatores. Reator tanque descontínuo (BSTR). Reator tanque de mistura contínuo (CSTR). Reator tubular de fluxo pistonado (PFR). Comparação de desempenho de reatores CSTR e PFR. Reatores CSTR em cascata. Associação mista de reatores em série: CSTR e PFR 5. ANÁLISE DE DADOS CINÉTICOS EM REATORES QUÍMICOS ISOTÉRMICOSBalanço de massa e coleta de dados em reatores ideais isotérmicos: batelada (BSTR), reator tanque de mistura contínuo (CSTR) e Reator tubular (PFR)'
$pPrograma.Range.Text = $programaPtText
$null = $pPrograma.Range.InsertParagraphAfter()
$pProgramaEn = $pPrograma.Next()
$programaEnText = '1. Introduction to KineticsTypes of Chemical Reactions.Rate law and its main parameters.Influence of temperature on reaction rate.Activation of chemical reactions.Arrhenius equation.Activation energy.Conversion.Concentration and its variation in a chemical transformation. 2. Reactions at Constant VolumeIrreversible reactions of first order.Irreversible reactions of second order.Irreversible reactions of third order.Irreversible reactions of any order. 3. Reactions at Variable VolumeConcepts.Volumetric conversion fraction.Reactions at variable volume of first and second order. 4. Ideal Models of Isothermal Chemical Reactors:Fundamental equations for reactor design.Batch reactor (BSTR).Continuous stirred-tank reactor (CSTR).Plug-flow reactor (PFR).Performance comparison of CSTR and PFR.Cascade CSTR reactors.Mixed association of reactors in series: CSTR and PFR. 5. Analysis of Kinetic Data in Isothermal Chemical ReactorsMass balance and data collection in ideal isothermal reactors:Batch reactor (BSTR).Continuous stirred-tank reactor (CSTR).Plug-flow reactor (PFR).'
$pProgramaEn.Range.InsertBefore($programaEnText)
$progEnStart = $pProgramaEn.Range.Start
$progEnRange = $d.Range($progEnStart, $progEnStart + $programaEnText.Length)
$progEnRange.Font.Italic = $true

# 6. Update "Método" evaluation text
$null = $d.Content.Find.Execute("Duas provas escritas (P1 e P2) e trabalhos relacionados à disciplina (TRAB).", $false, $false, $false, $false, $false, $true, 1, $false, 'Duas provas escritas (P1 e P2) e eventuais trabalhos relacionados à disciplina', 2)

# 7. Update "Critério" evaluation text
$null = $d.Content.Find.Execute("Média da Primeira Avaliação = (I)  Prova P1=30%; (II)  Prova P2=60% e (III)  Trabalhos =10%", $false, $false, $false, $false, $false, $true, 1, $false, 'Média da Primeira Avaliação (N) = 50% P1 + 50% P2.Obs: fica a critério de cada docente a inserção de trabalhos no decorrer do curso, bem como a alteração do peso de cada prova em decorrência dos mesmos.', 2)

# 8. Update "Norma de recuperação" evaluation text
$null = $d.Content.Find.Execute("Será a média aritmética da nota do aluno na primeira avaliação e da nota do aluo numa prova escrita na recuperação.", $false, $false, $false, $false, $false, $true, 1, $false, 'Média Final = (N + Prova Recuperação)/2', 2)

# 9. Rewrite "Bibliografia" body text (merge into one run with updated list)
$pBiblio = $null
foreach ($p in $d.Paragraphs) {
  if ($p.Range.Text.StartsWith("FOGLER")) { $pBiblio = $p }
}
$biblioText = '1- FOGLER, H.S. Elementos de engenharia das reações químicas. 3.ed. Rio de Janeiro: LTC Editora, 2009.2- LEVENSPIEL, O. Engenharia Das Reações Químicas, E ed (Blucher, São Paulo, 2000)3- VAN SANTEN, R.A.; Niemantsverdriet, J.W. Chemical kinetics and catalysis. New York: Plenum Press, 1995.4- Missen, R.W.; Mims, C.A.; Saville, B.A. Introduction to chemical reaction engineering and kinetics. New York: J. Wiley, 1999.5- Rothenberg, G. Catalysis: concepts and green applications. Weinheim: Wiley-VCH, 2008 Chichester.6- DENISOV, E.T.; Sarkisov, O.M.; Likhtenshtein, G.I. Chemical kinetics: fundamentals and new developments. Amsterdam: Elsevier, 2003.7- Hagen, J. Industrial catalysis: a practical approach. Weinheim: Wiley-VCH, 2006.8- Salmi, T.O.; Mikkola, J.; Warna, J.P. Chemical reaction engineering and reactor technology. Boca Raton: CRC Press/Taylor & Francis, 2011.9- Mortimer, M.; Taylor, P.G. Chemical kinetics and mechanism. Cambridge: Royal Society of Chemistry, 2002.10- FROMENT, G.F.; BISCHOFF, K.B. Chemical reactor analysis and design. 2nd. Ed. New York: John Wiley & Sons, 1990.11- HILL, C.G. An Introduction to chemical engineering kinetics and reactor design. New York: John Wiley&Sons, 1977.12- SMITH, J.M. Chemical engineering kinetics. 3rd. ed New York: McGraw-Hill,1981.13- DENBIGH, K.; TURNER, R. Introduction to chemical Reaction Design. Cambridge: Cambridge University Press, 1970.14 - AGUIAR, L. G. Problemas de cinética e reatores químicos. Curitiba: Appris Editora, 2023.'
$pBiblio.Range.Text = $biblioText

Write-Host "edit complete"
